# Add Modbus test cases to the jinzu connector test-data workbook.
#
# 1. Append a new worksheet "queryModbusEntity" after the last existing
#    sheet (getCacheKeyAndValue), populate its header/data rows (reusing
#    the same header layout + cell styles as the first sheet), and
#    2. Tweak the first sheet's view state (scroll/selection) to match the
#    authored change.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New sheet "queryModbusEntity"
# ---------------------------------------------------------------------
$headerSrc = $wb.Worksheets.Item(1)          # getConceptModelDataByCondition
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "queryModbusEntity"

# Header row (same columns/order as the other test sheets)
$ws.Range("A1").Value = "test-id"
$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "condition"
$ws.Range("D1").Value = "domainName"
$ws.Range("E1").Value = "fields"
$ws.Range("F1").Value = "name"
$ws.Range("G1").Value = "order"
$ws.Range("H1").Value = "pageIndex"
$ws.Range("I1").Value = "pageSize"
$ws.Range("J1").Value = "timeout"
$ws.Range("K1").Value = "rspStatus"
$ws.Range("L1").Value = "rspCode"
$ws.Range("M1").Value = "rspMessage"

# Data row - Modbus query test case
$ws.Range("A2").Value = "jinzu-modbus-query-var1"
$ws.Range("B2").Value = "good request, data retrieved"
$ws.Range("C2").Value = "stime > '`$start_time' and stime <  '`$end_time'"
$ws.Range("F2").Value = "ModbusTestEntity"

# Match formatting (header fill/border style + data-row border style) used
# by the other sheets.
$headerSrc.Range("A1:M1").Copy()
$ws.Range("A1:M1").PasteSpecial(-4122)   # xlPasteFormats
$headerSrc.Range("A2:M2").Copy()
$ws.Range("A2:M2").PasteSpecial(-4122)   # xlPasteFormats

# Column widths
$ws.Columns.Item(1).ColumnWidth = 36.44140625
$ws.Columns.Item(2).ColumnWidth = 29.88671875
$ws.Columns.Item(3).ColumnWidth = 62.33203125
$ws.Columns.Item(6).ColumnWidth = 19.77734375
$ws.Columns.Item(8).ColumnWidth = 24.88671875
$ws.Columns.Item(13).ColumnWidth = 22.5546875

$ws.Range("A11").Select()

# ---------------------------------------------------------------------
# 2. Update first sheet's view (scrolled further down, wider selection)
# ---------------------------------------------------------------------
$headerSrc.Range("A32").Select()
$excel.ActiveWindow.FreezePanes = $true
$headerSrc.Range("A1:XFD3").Select()
